$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column H: a "Save" indicator column, mirroring the header style of
# the existing header row (G1's format: bold, centered, bordered).
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("H1").Value = "Save"
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 1
$ws.Range("H4").Value = 1
$ws.Range("H5").Value = 0
